$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 499.5
$ws.Range("I107").Value = 499.5
$ws.Range("K107").Value = 499.5
$ws.Range("M107").Value = 1420.5
$ws.Range("H129").Value = 1706.25
$ws.Range("I129").Value = 1568
$ws.Range("K129").Value = 4704
$ws.Range("M129").Value = 296
$ws.Range("H135").Value = 446.7143
$ws.Range("I135").Value = 385.6
$ws.Range("J135").Value = 599.5
$ws.Range("K135").Value = 3470.4
$ws.Range("L135").Value = 5395.5
$ws.Range("M135").Value = -935.4000000000001
$ws.Range("N135").Value = -10465.5
$ws.Range("H138").Value = 5377.5386
$ws.Range("I138").Value = 2585.3333
$ws.Range("J138").Value = 5741.7393
$ws.Range("K138").Value = 7755.999899999999
$ws.Range("L138").Value = 17225.2179
$ws.Range("M138").Value = -2615.999899999999
$ws.Range("N138").Value = -27505.2179

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 900
$ws.Range("I2").Value = 900
$ws.Range("K2").Value = 900
$ws.Range("M2").Value = -787
$ws.Range("H10").Value = 1995
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H32").Value = 4526.56
$ws.Range("I32").Value = 4526.56
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4526.56
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4239.56
$ws.Range("H74").Value = 1804.5
$ws.Range("I74").Value = 1804.5
$ws.Range("K74").Value = 1804.5
$ws.Range("M74").Value = -930.5
$ws.Range("H77").Value = 1804.5
$ws.Range("I77").Value = 1804.5
$ws.Range("K77").Value = 9022.5
$ws.Range("M77").Value = -4654.5
$ws.Range("H102").Value = 2897.75
$ws.Range("J102").Value = 2898.5
$ws.Range("L102").Value = 2898.5
$ws.Range("N102").Value = -6142.5
$ws.Range("H116").Value = 900
$ws.Range("I116").Value = 900
$ws.Range("K116").Value = 900
$ws.Range("M116").Value = 1394
$ws.Range("H122").Value = 6766
$ws.Range("I122").Value = 6766
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 20298
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -17848
$ws.Range("H132").Value = 2750.8
$ws.Range("I132").Value = 1870.5385
$ws.Range("K132").Value = 5611.6155
$ws.Range("M132").Value = -3081.6155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 900
$ws.Range("I3").Value = 900
$ws.Range("K3").Value = 900
$ws.Range("M3").Value = -786
$ws.Range("H22").Value = 906.625
$ws.Range("I22").Value = 906.625
$ws.Range("K22").Value = 906.625
$ws.Range("M22").Value = -733.625
$ws.Range("H86").Value = 2847.6667
$ws.Range("I86").Value = 2232.7144
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 2232.7144
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1109.7144
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 2847.6667
$ws.Range("I89").Value = 2232.7144
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 11163.572
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -5547.572
$ws.Range("N89").Value = -36232
$ws.Range("H99").Value = 1098
$ws.Range("I99").Value = 1098
$ws.Range("K99").Value = 1098
$ws.Range("M99").Value = 400
$ws.Range("H107").Value = 800
$ws.Range("I107").Value = 800
$ws.Range("K107").Value = 800
$ws.Range("M107").Value = 1120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").ClearContents()
$ws.Range("N45").Value = 0
$ws.Range("H68").Value = 69795
$ws.Range("J68").Value = 69795
$ws.Range("L68").Value = 69795
$ws.Range("N68").Value = -71293
$ws.Range("H71").Value = 69795
$ws.Range("J71").Value = 69795
$ws.Range("L71").Value = 209385
$ws.Range("N71").Value = -216873
$ws.Range("H86").Value = 9781.75
$ws.Range("I86").Value = 8357
$ws.Range("J86").Value = 10799.429
$ws.Range("K86").Value = 8357
$ws.Range("L86").Value = 10799.429
$ws.Range("M86").Value = -7234
$ws.Range("N86").Value = -13045.429
$ws.Range("H89").Value = 9781.75
$ws.Range("I89").Value = 8357
$ws.Range("J89").Value = 10799.429
$ws.Range("K89").Value = 41785
$ws.Range("L89").Value = 53997.145
$ws.Range("M89").Value = -36169
$ws.Range("N89").Value = -65229.145
$ws.Range("H132").Value = 1692.3334
$ws.Range("J132").Value = 3257.2
$ws.Range("L132").Value = 9771.599999999999
$ws.Range("N132").Value = -14831.6
$ws.Range("H141").Value = 99992.5
$ws.Range("J141").Value = 99992.5
$ws.Range("L141").Value = 99992.5
$ws.Range("N141").Value = -110352.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 3126
$ws.Range("J93").Value = 3126
$ws.Range("L93").Value = 9378
$ws.Range("N93").Value = -13122

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 10000
$ws.Range("K29").Value = 10000
$ws.Range("M29").Value = -9710
$ws.Range("H122").Value = 9064.571
$ws.Range("I122").Value = 7454.154
$ws.Range("K122").Value = 22362.462
$ws.Range("M122").Value = -19912.462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2070.1428
$ws.Range("I22").Value = 1998.25
$ws.Range("K22").Value = 1998.25
$ws.Range("M22").Value = -1703.25
$ws.Range("H27").Value = 2070.1428
$ws.Range("I27").Value = 1998.25
$ws.Range("K27").Value = 1998.25
$ws.Range("M27").Value = -1891.25
$ws.Range("H40").Value = 1499.5
$ws.Range("I40").Value = 1499.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1499.5
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1363.5
$ws.Range("H46").Value = 3349.75
$ws.Range("I46").Value = 1699.5
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 1699.5
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -1511.5
$ws.Range("N46").Value = -5376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2449.5
$ws.Range("I126").Value = 2449.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7348.5
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -4878.5
$ws.Range("H136").Value = 1248.1052
$ws.Range("I136").Value = 1133.3846
$ws.Range("J136").Value = 1496.6666
$ws.Range("K136").Value = 3400.1538
$ws.Range("L136").Value = 4489.9998
$ws.Range("M136").Value = -850.1538
$ws.Range("N136").Value = -9589.9998
